$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Properties" table (B6:C16) is being restructured: the separate
# propertyNum/streetNum/streetName/postalCode/city/type fields collapse into
# a single "address" field, a new "sqFt" field is introduced, and the
# now-unused trailing rows are cleared.
$ws.Range("B7").Value = "address"
$ws.Range("C7").Value = "char"

$ws.Range("B8").Value = "city"
$ws.Range("C8").Value = "varChar"

$ws.Range("B9").Value = "postalCode"
$ws.Range("C9").Value = "varChar"

$ws.Range("B10").Value = "RoomCnt"
$ws.Range("C10").Value = "tinyInt"

$ws.Range("B11").Value = "BathroomCnt"
$ws.Range("C11").Value = "tinyInt"

$ws.Range("B12").Value = "description"
$ws.Range("C12").Value = "varChar"

$ws.Range("B13").Value = "sqFt"
$ws.Range("C13").Value = "int"

$ws.Range("B14").Value = "price"
$ws.Range("C14").Value = "int"

# Rows 15 and 16 no longer hold Properties-table entries.
$ws.Range("B15:C16").ClearContents()

# Selection moved from B17 to C12.
$ws.Range("C12").Select()
